## Generate Report for Handoff
## Adds a new handoff row (98fe604c-a387-475a-9d4e-04e169db49e5.md) to the
## Overview / zh-cn / de-de localization-status sheets, mirroring the
## existing 6b18567d-... row.

$wb = $excel.ActiveWorkbook

$fileBase = "98fe604c-a387-475a-9d4e-04e169db49e5"
$fileMd   = "$fileBase.md"
$zhXlf    = "$fileBase.42bc41eff598040fe763e49d5cc08f9d5ee0069b.zh-cn.xlf"
$deXlf    = "$fileBase.42bc41eff598040fe763e49d5cc08f9d5ee0069b.de-de.xlf"
$ghUrl    = "https://github.com/OpenLocalizationTestOrg/oltest/blob/8166c3f78769bc3dfad65bbf07c67c2dabd75889/e2e/$fileMd"
$dateFmt  = "yyyy-mm-dd HH:mm:ss"

## ---------------------------------------------------------------
## Overview sheet (sheet1) -> new row 3
## ---------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A3").Value = $fileMd
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("D3").Value = "'"
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-13 02:48:06"
$wsOverview.Range("G3").NumberFormat = $dateFmt

$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $ghUrl, $null, $null, "e2e\$fileMd") | Out-Null

$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:G3"))

## ---------------------------------------------------------------
## zh-cn sheet (sheet2) -> new row 3
## ---------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("D3").Value = "e2e"
$wsZhCn.Range("E3").Value = "ht"
$wsZhCn.Range("F3").Value = "'False"
$wsZhCn.Range("G3").Value = $zhXlf
$wsZhCn.Range("H3").Value = "2016-08-13 02:47:55"
$wsZhCn.Range("H3").NumberFormat = $dateFmt
$wsZhCn.Range("I3").Value = "'"
$wsZhCn.Range("J3").Value = "'"
$wsZhCn.Range("K3").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("L3").Value = "'"
$wsZhCn.Range("M3").Value = "'True"
$wsZhCn.Range("N3").Value = "'"
$wsZhCn.Range("O3").Value = "'False"
$wsZhCn.Range("P3").Value = "'"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $ghUrl, $null, $null, $fileMd) | Out-Null

$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.Resize($wsZhCn.Range("A1:P3"))

## ---------------------------------------------------------------
## de-de sheet (sheet3) -> new row 3
## ---------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("B3").Value = ".md"
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("D3").Value = "e2e"
$wsDeDe.Range("E3").Value = "ht"
$wsDeDe.Range("F3").Value = "'False"
$wsDeDe.Range("G3").Value = $deXlf
$wsDeDe.Range("H3").Value = "2016-08-13 02:48:06"
$wsDeDe.Range("H3").NumberFormat = $dateFmt
$wsDeDe.Range("I3").Value = "'"
$wsDeDe.Range("J3").Value = "'"
$wsDeDe.Range("K3").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("L3").Value = "'"
$wsDeDe.Range("M3").Value = "'True"
$wsDeDe.Range("N3").Value = "'"
$wsDeDe.Range("O3").Value = "'False"
$wsDeDe.Range("P3").Value = "'"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $ghUrl, $null, $null, $fileMd) | Out-Null

$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.Resize($wsDeDe.Range("A1:P3"))
